# Updates crypto price/volume table per latest scrape (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.459.55'
$ws.Range('E2').Value = '  -1.48%  '
$ws.Range('D3').Value = '2.533.99'
$ws.Range('E3').Value = '  -1.21%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''311.48'
$ws.Range('E5').Value = '  -1.66%  '
$ws.Range('D6').Value = '''98.63'
$ws.Range('E6').Value = '  +2.02%  '
$ws.Range('D7').Value = '''0.567'
$ws.Range('E7').Value = '  -1.43%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('E9').Value = '  -2.48%  '
$ws.Range('D10').Value = '''35.67'
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('E11').Value = '  -1.37%  '
$ws.Range('E12').Value = '  -1.61%  '
$ws.Range('E13').Value = '  -0.19%  '
$ws.Range('D14').Value = '2.923.03'
$ws.Range('E14').Value = '  -1.20%  '
$ws.Range('D15').Value = '''15.75'
$ws.Range('E15').Value = '  +4.83%  '
$ws.Range('D16').Value = '2.486.37'
$ws.Range('E16').Value = '  -3.42%  '
$ws.Range('D17').Value = '''0.828'
$ws.Range('E17').Value = '  -2.13%  '
$ws.Range('D18').Value = '42.468.29'
$ws.Range('E18').Value = '  -1.52%  '
$ws.Range('D19').Value = '''6.79'
$ws.Range('E19').Value = '  -1.08%  '
$ws.Range('D20').Value = '0.0₃0949'
$ws.Range('E20').Value = '  -1.36%  '
$ws.Range('D21').Value = '''12.22'
$ws.Range('E21').Value = '  -3.07%  '
$ws.Range('D22').Value = '''69.04'
$ws.Range('E22').Value = '  -0.72%  '
$ws.Range('D23').Value = '''243.73'
$ws.Range('E23').Value = '  -3.66%  '
$ws.Range('E24').Value = '  -1.52%  '
$ws.Range('D25').Value = '''2.05'
$ws.Range('E25').Value = '  -0.75%  '
$ws.Range('E27').Value = '  -3.03%  '
$ws.Range('D28').Value = '''2.33'
$ws.Range('E28').Value = '  -4.32%  '
$ws.Range('D29').Value = '''39.29'
$ws.Range('E29').Value = '  -2.25%  '
$ws.Range('D30').Value = '''10.13'
$ws.Range('E30').Value = '  -0.95%  '
$ws.Range('D31').Value = '''157.76'
$ws.Range('E31').Value = '  +2.43%  '
$ws.Range('E32').Value = '  -2.20%  '
$ws.Range('E33').Value = '  +15.46%  '
$ws.Range('D34').Value = '''0.0796'
$ws.Range('E34').Value = '  -1.51%  '
$ws.Range('E35').Value = '  -3.20%  '
$ws.Range('E36').Value = '  -5.15%  '
$ws.Range('E37').Value = '  -7.72%  '
$ws.Range('D38').Value = '''18.11'
$ws.Range('E38').Value = '  -4.79%  '
$ws.Range('E39').Value = '  -0.67%  '
$ws.Range('E40').Value = '  -0.29%  '
$ws.Range('D41').Value = '''4.27'
$ws.Range('E41').Value = '  +9.32%  '
$ws.Range('D42').Value = '''21.63'
$ws.Range('E42').Value = '  -4.07%  '
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('D44').Value = '''3.29'
$ws.Range('E44').Value = '  +0.68%  '
$ws.Range('E45').Value = '  -2.55%  '
$ws.Range('D46').Value = '1.963.53'
$ws.Range('E46').Value = '  -1.65%  '
$ws.Range('D47').Value = '''8.94'
$ws.Range('E47').Value = '  -0.91%  '
$ws.Range('D48').Value = '2.777.10'
$ws.Range('E48').Value = '  -1.37%  '
$ws.Range('D49').Value = '''81.00'
$ws.Range('E49').Value = '  -3.49%  '
$ws.Range('E50').Value = '  -0.82%  '
$ws.Range('B51').Value = 'SEI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range('D51').Value = '''0.844'
$ws.Range('E51').Value = '  +7.82%  '
